# 5.2.1.1b: add a new "2023" column (R) to the table, mirroring column Q
# (2022) for layout/styling, then fill in the handful of cells that carry
# genuinely new data. All other 2023 cells reuse the same "…" (no data)
# placeholder already used throughout column Q for rows without figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy Q3:Q25 -> R3:R25 so the new column inherits the exact same styles
# (borders, number formats, "…" placeholders, etc.) as column Q.
$ws.Range("Q3:Q25").Copy($ws.Range("R3:R25"))

# Header year
$ws.Range("R4").Value = 2023

# New 2023 figures (rows with real data)
$ws.Range("R5").Value = 11357
$ws.Range("R7").Value = 11002
$ws.Range("R8").Value = 355
